$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores Price (column D) and Volume(1h) (column E) as
# literal text (e.g. "236.22", "0.611", thousands-dotted "40.994.33", and
# padded percentages like "  -2.42%  "). Values that look like a plain
# decimal number would otherwise be auto-converted to a Number by the
# Value setter, so force those specific cells to Text format first.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '41.017.56'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '2.170.61'
$ws.Range("E3").Value = '  -2.36%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '236.22'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  -2.69%  '
$ws.Range("D7").Value = '69.46'
$ws.Range("E7").Value = '  -5.76%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -6.83%  '
$ws.Range("D10").Value = '39.77'
$ws.Range("E10").Value = '  -8.14%  '
$ws.Range("D11").Value = '0.0924'
$ws.Range("E11").Value = '  -3.79%  '
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("D13").Value = '6.73'
$ws.Range("E13").Value = '  -5.44%  '
$ws.Range("D14").Value = '2.495.60'
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = '13.79'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").Value = '0.808'
$ws.Range("E16").Value = '  -4.36%  '
$ws.Range("D17").Value = '2.169.41'
$ws.Range("E17").Value = '  -3.09%  '
$ws.Range("D18").Value = '40.847.64'
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("E19").Value = '  -7.61%  '
$ws.Range("D20").Value = '70.30'
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").Value = '5.92'
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("D22").Value = '9.61'
$ws.Range("E22").Value = '  -5.48%  '
$ws.Range("D23").Value = '225.02'
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  -7.90%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '10.89'
$ws.Range("E26").Value = '  -6.16%  '
$ws.Range("D27").Value = '3.53'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").Value = '165.96'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").Value = '19.81'
$ws.Range("E31").Value = '  -4.02%  '
$ws.Range("D32").Value = '30.76'
$ws.Range("E32").Value = '  +4.30%  '
$ws.Range("D33").Value = '0.0772'
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  -8.54%  '
$ws.Range("D35").Value = '0.120'
$ws.Range("E35").Value = '  -3.50%  '
$ws.Range("E36").Value = '  -9.49%  '
$ws.Range("D37").Value = '4.13'
$ws.Range("E37").Value = '  -4.31%  '
$ws.Range("D38").Value = '0.0286'
$ws.Range("E38").Value = '  -5.07%  '
$ws.Range("D39").Value = '12.37'
$ws.Range("E39").Value = '  -4.98%  '
$ws.Range("D40").Value = '2.05'
$ws.Range("E40").Value = '  -4.08%  '
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("D42").Value = '59.68'
$ws.Range("E42").Value = '  -7.71%  '
$ws.Range("E43").Value = '  -5.60%  '
$ws.Range("D44").Value = '8.26'
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("D45").Value = '0.0971'
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("D46").Value = '98.74'
$ws.Range("E46").Value = '  -5.70%  '
$ws.Range("E47").Value = '  -3.36%  '
$ws.Range("D48").Value = '1.12'
$ws.Range("E48").Value = '  -3.36%  '
$ws.Range("E49").Value = '  -7.92%  '
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").Value = '2.373.76'
$ws.Range("E51").Value = '  -2.38%  '
